$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-29). The workbook was refreshed a day later, so bump each of
# those date values by one day (45575 -> 45576).
$ws.Range("C2:C29").Value = 45576
